# Generate Report for Handoff
# Update "Latest Handoff Datetime" (column D) for rows that were just
# re-handed-off, on both the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$rowsToUpdate = @(4, 6, 7, 8, 9, 10)

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rowsToUpdate) {
    $ws_zhcn.Range("D$r").Value = "2016-03-04 06:55:25"
}

$ws_dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rowsToUpdate) {
    $ws_dede.Range("D$r").Value = "2016-03-04 06:55:41"
}
